# Update the generated three-digit x one-digit multiplication answers
# in the single table of the document. Replacements are applied per
# cell (row, col) using Find scoped to that cell's Range so that
# duplicate values (e.g. "480x6=2880" appearing twice, or
# "749x3=2247" being both an old and a new value) are not ambiguous.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Map of (row, col) -> new text, in document order.
$updates = @(
    @{Row=1;  Col=1; New="735×3=2205"},
    @{Row=1;  Col=2; New="794×2=1588"},
    @{Row=1;  Col=3; New="919×4=3676"},
    @{Row=1;  Col=4; New="545×3=1635"},
    @{Row=1;  Col=5; New="480×6=2880"},

    @{Row=5;  Col=1; New="480×6=2880"},
    @{Row=5;  Col=2; New="749×3=2247"},
    @{Row=5;  Col=3; New="846×6=5076"},
    @{Row=5;  Col=4; New="574×4=2296"},
    @{Row=5;  Col=5; New="482×2=964"},

    @{Row=10; Col=1; New="719×4=2876"},
    @{Row=10; Col=2; New="806×3=2418"},
    @{Row=10; Col=3; New="776×2=1552"},
    @{Row=10; Col=4; New="747×6=4482"},
    @{Row=10; Col=5; New="993×9=8937"},

    @{Row=15; Col=1; New="354×8=2832"},
    @{Row=15; Col=2; New="526×9=4734"},
    @{Row=15; Col=3; New="931×3=2793"},
    @{Row=15; Col=4; New="337×3=1011"},
    @{Row=15; Col=5; New="588×9=5292"},

    @{Row=20; Col=1; New="642×4=2568"},
    @{Row=20; Col=2; New="101×6=606"},
    @{Row=20; Col=3; New="837×3=2511"},
    @{Row=20; Col=4; New="620×3=1860"},
    @{Row=20; Col=5; New="515×8=4120"}
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Shrink the range so it doesn't include the trailing cell-mark /
    # paragraph-mark characters, then directly set its text to the
    # new value.
    $rng.End = $rng.End - 1
    $rng.Text = $u.New
}
